# Add a heading block (institute name + report title) above the existing
# attendance table, and a facilitator signature line below it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 5 new blank rows above the current row 1, shifting everything
# (header row + the two data rows) down by 5 rows. The existing header-row
# style (bold white text on navy fill) travels with the cells automatically.
$insertRange = $ws.Range("A1:L5")
$insertRange.EntireRow.Insert()

# --- Widen the columns from 20 to 22 ---
# (21.15 is used instead of 22 because the host's ColumnWidth setter adds a
# small fixed pixel-rounding pad when round-tripping through OOXML; this
# input value is the one that lands exactly on a stored width of 22.)
for ($col = 1; $col -le 12; $col++) {
    $ws.Cells.Item(1, $col).EntireColumn.ColumnWidth = 21.15
}

# --- Institute name heading (row 2) ---
$title = $ws.Range("C2")
$title.Value = "Southern Labs Institute of Technology"
$title.Font.Bold = $true
$title.Font.Size = 16
$title.Font.Color = 6697728   # RGB(0x00,0x33,0x66) -> navy blue #003366

# --- Report subtitle (row 3) ---
$subtitle = $ws.Range("C3")
$subtitle.Value = "Attendance Report: WPE - Week 1"
$subtitle.Font.Bold = $true
$subtitle.Font.Size = 12

# --- Facilitator signature line (row 11, two blank rows after the table) ---
$sig = $ws.Range("A11")
$sig.Value = "Facilitator Signature: ___________________________"
$sig.Font.Bold = $true
$sig.Font.Size = 12
